$d = $word.ActiveDocument

$find = "Kampagnendaten Sternbild Herkules 2022: 13.-22. Juni, 12.-21. Juli, 10.-19. August"
$replace = "Kampagnendaten 2022 für das Sternbild Sternbild Herkules: 13.-22. Juni, 12.-21. Juli, 10.-19. August"

$rng = $d.Content
$rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
